# Add the Axure prototype summary notes to the "Fred Lei" sheet,
# and make it the active sheet (was "Aniket").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fred Lei")

$ws.Range("A12").Value = "the  main page protype. "
$ws.Range("B12").Value = "A.  5th,May   Fred&Yako  meeting minites "

$ws.Range("B14").Value = "1.   using  Axure to draw the prototype"
$ws.Range("B15").Value = "2.   common features(login, sign up, setting,...)   (owner: Fred Lei)"
$ws.Range("B16").Value = "3.   main page:  (owner: Yako)"
$ws.Range("B17").Value = "(top:  background (nice picture) and where to go function."
$ws.Range("B18").Value = " next: Category by mood with small icon :  wedding, luxury,  sport, education, bar-life, peace"
$ws.Range("B19").Value = "Next:  hot products list(nice picture with short words)  (can be clicked and jump to detail page)"
$ws.Range("B20").Value = "Next:  customer blog (video, image, type of mood, address, time, (copy from douyin))"
$ws.Range("B21").Value = "main function bars:   my Trips,  create trip, post blog, write review, me(setting)"
$ws.Range("B22").Value = " )"

$ws.Activate()
$ws.Range("B29").Select()
